# Apply the "PLO Stats-this session" -> "NL Stats-this session" rename,
# refresh the latest-session stats rows, drop the now-empty last data row
# (Xavier), and repoint every chart series formula so cached cat/val
# ranges and sheet names stay in sync with the new 3-row data block.

$wb = $excel.ActiveWorkbook

$combinedName = "combined Stats-this session"
$oldSessionName = "PLO Stats-this session"
$newSessionName = "NL Stats-this session"

# ---------------------------------------------------------------------
# 1. Rename the "PLO Stats-this session" sheet to "NL Stats-this session"
# ---------------------------------------------------------------------
$wb.Worksheets.Item($oldSessionName).Name = $newSessionName

# ---------------------------------------------------------------------
# 2. Refresh the data rows on both sheets (they mirror each other)
# ---------------------------------------------------------------------

# Row 2: Fish
$row2 = @{
    A = "Fish"; B = 20; C = 53.59; D = 33.59; E = 0
    F = 0.679; G = 0; H = 0; I = 0.118; J = 0.268; K = 0.19
    L = 0.19; M = 0; N = 0
    O = 112.17; P = 66.58; Q = 168; R = 0.711
}

# Row 3: Raymond (was Cedric)
$row3 = @{
    A = "Raymond"; B = 51; C = 24.39; D = -26.61; E = 0
    F = 0.466; G = 0.31; H = 0.0057; I = 0.454; J = 0.126; K = 0.046
    L = 3.68; M = 26; N = 56
    O = 41.7; P = 58.54; Q = 174; R = 0.364
}

# Row 4: Scott
$row4 = @{
    A = "Scott"; B = 20; C = 13.02; D = -6.98; E = 0
    F = 0.527; G = 0.24; H = 0.02; I = 0.513; J = 0.22; K = 0.08
    L = 3.16; M = 16; N = 45
    O = 43.65; P = 64.11; Q = 150; R = 0.364
}

foreach ($sheetName in @($combinedName, $newSessionName)) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($col in $row2.Keys) { $ws.Range($col + "2").Value = $row2[$col] }
    foreach ($col in $row3.Keys) { $ws.Range($col + "3").Value = $row3[$col] }
    foreach ($col in $row4.Keys) { $ws.Range($col + "4").Value = $row4[$col] }

    # T3 holds a date-shaped literal string ("07/05/21"), not a real date
    # serial - force text formatting while assigning so it doesn't get
    # auto-parsed into a date, then drop back to the default style so the
    # cell ends up unstyled, same as the rest of the sheet.
    $ws.Range("T3").NumberFormat = "@"
    $ws.Range("T3").Value = "07/05/21"
    $ws.Range("T3").Style = "Normal"

    # -------------------------------------------------------------
    # 3. Drop row 5 (Xavier) entirely - only 3 data rows remain now
    # -------------------------------------------------------------
    $ws.Range("A5").EntireRow.Delete()
}

# ---------------------------------------------------------------------
# 4. Repoint every chart series: cached cat/val ranges go from row 2-5
#    to row 2-4, and any chart still pointing at the old sheet name picks
#    up the new one. The series "name" reference (e.g. $M$1) is edited
#    via .Name using a relative A1 ref so it round-trips the same way the
#    original file stored it (no forced $ absolute markers), while cat/val
#    are edited via .XValues/.Values which already used absolute refs.
# ---------------------------------------------------------------------
function Update-ChartSeriesRanges($ws, $newSheetName) {
    foreach ($co in $ws.ChartObjects()) {
        $chart = $co.Chart
        $series = $chart.SeriesCollection()
        for ($i = 1; $i -le $series.Count; $i++) {
            $s = $series.Item($i)
            $f = $s.Formula
            $pattern = "^=SERIES\('[^']*'!(\`$?[A-Z]+\`$?[0-9]+),'[^']*'!(\`$[A-Z]+\`$[0-9]+:\`$[A-Z]+\`$)5,'[^']*'!(\`$[A-Z]+\`$[0-9]+:\`$[A-Z]+\`$)5,(.*)\)$"
            if ($f -match $pattern) {
                $nameCellRel = $matches[1] -replace '\$', ''
                $catPrefix = $matches[2]
                $valPrefix = $matches[3]
                $s.Name = "='" + $newSheetName + "'!" + $nameCellRel
                $s.XValues = "='" + $newSheetName + "'!" + $catPrefix + "4"
                $s.Values = "='" + $newSheetName + "'!" + $valPrefix + "4"
            }
        }
    }
}

Update-ChartSeriesRanges $wb.Worksheets.Item($combinedName) $combinedName
Update-ChartSeriesRanges $wb.Worksheets.Item($newSessionName) $newSessionName
